$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add a new paragraph "P – rilascia Rescue Pack" right after the paragraph
#    "D - richiesta di informazione di un determinato sensore", carrying the
#    _GoBack bookmark.
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("di un determinato sensore", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find anchor paragraph for Rescue Pack insertion"
}
# Collapse to the end of the match (just before the paragraph mark) and
# insert a brand new paragraph there.
$insertionPoint = $d.Range($rng1.End, $rng1.End)

$rescuePackPkg = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:rPr><w:i/></w:rPr><w:t>P</w:t></w:r><w:r><w:t xml:space="preserve"> – rilascia Rescue Pack</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertionPoint.InsertXML($rescuePackPkg)

# ---------------------------------------------------------------------------
# 2) Bump the z-index of the floating VML picture (posizione_sensori) from
#    -251657216 to -251658752, leaving everything else about the shape as-is.
# ---------------------------------------------------------------------------
$shapeParaIndex = -1
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.WordOpenXML.Contains("z-index")) {
        $shapeParaIndex = $i
        break
    }
}
if ($shapeParaIndex -eq -1) {
    throw "Could not find paragraph containing the floating picture"
}
$shapePara = $d.Paragraphs($shapeParaIndex)
$shapeRng = $d.Range($shapePara.Range.Start, $shapePara.Range.End)

$shapePkg = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w:rsidR="00A25C49" w:rsidRDefault="005229E3" w:rsidP="00A25C49"><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:pict><v:shapetype id="_x0000_t75" coordsize="21600,21600" o:spt="75" o:preferrelative="t" path="m@4@5l@4@11@9@11@9@5xe" filled="f" stroked="f"><v:stroke joinstyle="miter"/><v:formulas><v:f eqn="if lineDrawn pixelLineWidth 0"/><v:f eqn="sum @0 1 0"/><v:f eqn="sum 0 0 @1"/><v:f eqn="prod @2 1 2"/><v:f eqn="prod @3 21600 pixelWidth"/><v:f eqn="prod @3 21600 pixelHeight"/><v:f eqn="sum @0 0 1"/><v:f eqn="prod @6 1 2"/><v:f eqn="prod @7 21600 pixelWidth"/><v:f eqn="sum @8 21600 0"/><v:f eqn="prod @7 21600 pixelHeight"/><v:f eqn="sum @10 21600 0"/></v:formulas><v:path o:extrusionok="f" gradientshapeok="t" o:connecttype="rect"/><o:lock v:ext="edit" aspectratio="t"/></v:shapetype><v:shape id="_x0000_s1026" type="#_x0000_t75" style="position:absolute;margin-left:222.65pt;margin-top:7.6pt;width:246.05pt;height:175.7pt;z-index:-251658752;mso-position-horizontal-relative:text;mso-position-vertical-relative:text;mso-width-relative:page;mso-height-relative:page"><v:imagedata r:id="rId4" o:title="posizione_sensori" croptop="7989f" cropbottom="15306f" cropleft="12788f" cropright="8484f"/></v:shape></w:pict></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$shapeRng.InsertXML($shapePkg)

# ---------------------------------------------------------------------------
# 3) Fix the "Temperatura" sensor row: merge the two runs that were split by
#    the stray _GoBack bookmark back into a single run / remove the bookmark.
#    A self-replace over the whole visible text forces Word to rebuild the
#    run (dropping the now-redundant bookmark) without touching anything
#    else in the cell.
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("Temperatura (20 – 40)", $true, $false, $false, $false, $false, $true, 1, $false, "Temperatura (20 – 40)", 2)
if (-not $found3) {
    throw "Could not find Temperatura sensor text"
}

Write-Output "done"
